# Update cryptocurrency price/volume figures per the Feb 17 2023 data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ensure the price (D) and volume (E) columns remain plain text so values like
# "0.02600" or "-2.85%" keep their exact original formatting instead of being
# reinterpreted by Excel as numbers/percentages.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '310.21'
$ws.Range("E2").Value = '-2.85%'
$ws.Range("D3").Value = '51.53'
$ws.Range("E3").Value = '5.63%'
$ws.Range("D4").Value = '5.128'
$ws.Range("E4").Value = '-2.78%'
$ws.Range("D5").Value = '0.07786'
$ws.Range("E5").Value = '-3.30%'
$ws.Range("D6").Value = '4.522'
$ws.Range("E6").Value = '-1.45%'
$ws.Range("D7").Value = '1.359'
$ws.Range("E7").Value = '1.71%'
$ws.Range("D8").Value = '1.593'
$ws.Range("E8").Value = '-3.87%'
$ws.Range("D9").Value = '0.1218'
$ws.Range("E9").Value = '-6.20%'
$ws.Range("E10").Value = '2.73%'
$ws.Range("D11").Value = '0.09722'
$ws.Range("E11").Value = '4.14%'
$ws.Range("D12").Value = '0.04723'
$ws.Range("E12").Value = '3.05%'
$ws.Range("D13").Value = '0.1046'
$ws.Range("E13").Value = '0.15%'
$ws.Range("D14").Value = '0.001259'
$ws.Range("E14").Value = '-5.57%'
$ws.Range("D15").Value = '0.005818'
$ws.Range("E15").Value = '-0.49%'
$ws.Range("E16").Value = '2,015.33%'
$ws.Range("E17").Value = '-0.06%'
$ws.Range("E18").Value = '-0.25%'
$ws.Range("D19").Value = '0.3471'
$ws.Range("E19").Value = '1.52%'
$ws.Range("D20").Value = '7.982'
$ws.Range("E20").Value = '-2.08%'
$ws.Range("D21").Value = '0.1372'
$ws.Range("E21").Value = '-2.09%'
$ws.Range("D22").Value = '0.3089'
$ws.Range("E22").Value = '-0.29%'
$ws.Range("D23").Value = '0.04159'
$ws.Range("E23").Value = '-0.23%'
$ws.Range("E24").Value = '-2.81%'
$ws.Range("D25").Value = '0.004039'
$ws.Range("E25").Value = '-4.96%'
$ws.Range("E26").Value = '-0.05%'
$ws.Range("D38").Value = '0.02600'
$ws.Range("E38").Value = '-3.69%'
$ws.Range("D39").Value = '0.05952'
$ws.Range("E39").Value = '4.42%'
$ws.Range("D40").Value = '0.01113'
$ws.Range("E40").Value = '76.24%'
$ws.Range("D41").Value = '0.008074'
$ws.Range("E41").Value = '1.26%'
$ws.Range("D42").Value = '0.1423'
$ws.Range("E42").Value = '-1.40%'
$ws.Range("E43").Value = '6.89%'
$ws.Range("D44").Value = '0.007730'
$ws.Range("E44").Value = '-1.34%'
$ws.Range("D45").Value = '0.3130'
$ws.Range("E45").Value = '-10.73%'
$ws.Range("D46").Value = '0.00007254'
$ws.Range("E46").Value = '5.15%'
$ws.Range("D47").Value = '0.00000000750'
$ws.Range("D48").Value = '0.05317'
$ws.Range("E48").Value = '-3.16%'
$ws.Range("D49").Value = '0.002619'
$ws.Range("E49").Value = '-34.53%'
$ws.Range("D50").Value = '0.00002100'
$ws.Range("D51").Value = '0.0002000'
